$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) keeps its text formatting so numeric-looking
# strings like "595.16" are not auto-converted into Excel numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '67.436.66'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").Value = '2.628.70'
$ws.Range("E3").Value = '  -1.51%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '595.16'
$ws.Range("E5").Value = '  -0.65%  '
$ws.Range("D6").Value = '167.86'
$ws.Range("E6").Value = '  +1.42%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '0.535'
$ws.Range("E8").Value = '  -2.31%  '
$ws.Range("D9").Value = '2.628.79'
$ws.Range("E9").Value = '  -1.48%  '
$ws.Range("D10").Value = '0.140'
$ws.Range("E10").Value = '  -1.92%  '
$ws.Range("E11").Value = '  +1.37%  '
$ws.Range("D12").Value = '0.366'
$ws.Range("E12").Value = '  +2.16%  '
$ws.Range("E13").Value = '  +0.50%  '
$ws.Range("D14").Value = '27.71'
$ws.Range("E14").Value = '  -0.19%  '
$ws.Range("D15").Value = '3.106.38'
$ws.Range("E15").Value = '  -1.68%  '
$ws.Range("D16").Value = '0.0000183'
$ws.Range("E16").Value = '  -0.49%  '
$ws.Range("D17").Value = '67.212.65'
$ws.Range("E17").Value = '  -0.24%  '
$ws.Range("D18").Value = '2.611.86'
$ws.Range("E18").Value = '  -2.28%  '
$ws.Range("D19").Value = '12.11'
$ws.Range("E19").Value = '  +3.46%  '
$ws.Range("D20").Value = '8.01'
$ws.Range("E20").Value = '  +5.15%  '
$ws.Range("D21").Value = '358.90'
$ws.Range("E21").Value = '  -1.30%  '
$ws.Range("E22").Value = '  -0.84%  '
$ws.Range("D23").Value = '4.69'
$ws.Range("E23").Value = '  -2.46%  '
$ws.Range("B24").Value = 'SuiNetwork'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D24").Value = '1.94'
$ws.Range("E24").Value = '  -4.73%  '
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("D26").Value = '10.29'
$ws.Range("E26").Value = '  +1.56%  '
$ws.Range("D27").Value = '69.83'
$ws.Range("E27").Value = '  -1.59%  '
$ws.Range("E29").Value = '  +0.11%  '
$ws.Range("D30").Value = '0.0000101'
$ws.Range("E30").Value = '  -1.21%  '
$ws.Range("D31").Value = '544.96'
$ws.Range("E31").Value = '  -1.84%  '
$ws.Range("D32").Value = '7.95'
$ws.Range("E32").Value = '  -0.23%  '
$ws.Range("E33").Value = '  -2.41%  '
$ws.Range("D34").Value = '1.91'
$ws.Range("E34").Value = '  -1.21%  '
$ws.Range("E35").Value = '  +5.06%  '
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("E37").Value = '  -2.92%  '
$ws.Range("D38").Value = '156.72'
$ws.Range("E38").Value = '  +1.40%  '
$ws.Range("D39").Value = '19.04'
$ws.Range("E39").Value = '  -2.50%  '
$ws.Range("D40").Value = '0.367'
$ws.Range("E40").Value = '  -1.57%  '
$ws.Range("D41").Value = '5.23'
$ws.Range("E41").Value = '  -1.26%  '
$ws.Range("E42").Value = '  -0.52%  '
$ws.Range("D43").Value = '18.21'
$ws.Range("E43").Value = '  +1.44%  '
$ws.Range("E44").Value = '  +0.04%  '
$ws.Range("D45").Value = '2.43'
$ws.Range("E45").Value = '  -3.55%  '
$ws.Range("D46").Value = '0.0₆0297'
$ws.Range("E46").Value = '  -0.29%  '
$ws.Range("D47").Value = '152.55'
$ws.Range("E47").Value = '  -0.25%  '
$ws.Range("D48").Value = '0.581'
$ws.Range("E48").Value = '  -1.32%  '
$ws.Range("D49").Value = '3.79'
$ws.Range("E49").Value = '  -0.92%  '
$ws.Range("E50").Value = '  -0.61%  '
$ws.Range("E51").Value = '  -1.06%  '
